# Update "想去人数" (F column) values across the workbook's sheets.
# Values were regenerated from a later scrape, so only column F numbers change.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 7930
$ws1.Range("F3").Value = 109
$ws1.Range("F4").Value = 86
$ws1.Range("F5").Value = 17636
$ws1.Range("F8").Value = 657
$ws1.Range("F14").Value = 167
$ws1.Range("F15").Value = 354
$ws1.Range("F17").Value = 307
$ws1.Range("F18").Value = 145
$ws1.Range("F19").Value = 401
$ws1.Range("F21").Value = 1104
$ws1.Range("F22").Value = 81
$ws1.Range("F23").Value = 661
$ws1.Range("F24").Value = 2247
$ws1.Range("F25").Value = 773
$ws1.Range("F26").Value = 61
$ws1.Range("F27").Value = 568
$ws1.Range("F29").Value = 623
$ws1.Range("F30").Value = 561

# Sheet "演出" (Performance)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 294
$ws2.Range("F3").Value = 68

# Sheet "本地生活" (Local Life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 490

# Sheet "全部类型" (All Types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 490
$ws4.Range("F3").Value = 7930
$ws4.Range("F4").Value = 109
$ws4.Range("F5").Value = 86
$ws4.Range("F6").Value = 294
$ws4.Range("F7").Value = 17642
$ws4.Range("F10").Value = 657
$ws4.Range("F12").Value = 68
$ws4.Range("F20").Value = 169
$ws4.Range("F21").Value = 354
$ws4.Range("F27").Value = 307
$ws4.Range("F28").Value = 145
$ws4.Range("F29").Value = 401
$ws4.Range("F31").Value = 1104
$ws4.Range("F32").Value = 81
$ws4.Range("F33").Value = 661
$ws4.Range("F34").Value = 2248
$ws4.Range("F35").Value = 773
$ws4.Range("F36").Value = 61
$ws4.Range("F37").Value = 568
$ws4.Range("F40").Value = 623
$ws4.Range("F41").Value = 561
